$d = $word.ActiveDocument

# 1) Replace the long heading paragraph text with the {header} placeholder.
$d.Paragraphs.Item(1).Range.Text = "{header}"

# 2) Replace the table header cell captions with their placeholders
#    (setting Range.Text keeps the existing run's <w:rPr/> intact, unlike Find/Replace).
$table = $d.Tables.Item(1)
$table.Cell(1, 1).Range.Text = "{col1}"
$table.Cell(1, 2).Range.Text = "{col2}"
$table.Cell(1, 3).Range.Text = "{col3}"
$table.Cell(1, 4).Range.Text = "{col4}"
$table.Cell(1, 5).Range.Text = "{col5}"

# 3) Fix the rounding on the last grid column / 5th cell width (1871 -> 1872 dxa).
$table.Columns.Item(5).Width = 93.6
